$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D, shifting "Tipo" (and its data) from D to E
$ws.Columns.Item(4).Insert()

# Set the new header and value in column D
$ws.Cells.Item(1, 4).Value = "MAE"
$ws.Cells.Item(2, 4).Value = 2.091210877959376

# Copy the style of the existing header cell (C1) onto the new header cell (D1)
$ws.Cells.Item(1, 3).Copy()
$ws.Cells.Item(1, 4).PasteSpecial(-4122)
$excel.CutCopyMode = $false
